$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 1316.9412
$ws.Range("I58").Value = 138.8
$ws.Range("J58").Value = 3000
$ws.Range("K58").Value = 416.4
$ws.Range("L58").Value = 9000
$ws.Range("M58").Value = -266.4
$ws.Range("N58").Value = -9300
$ws.Range("H137").Value = 941.4815
$ws.Range("I137").Value = 810.4545000000001
$ws.Range("K137").Value = 2431.3635
$ws.Range("M137").Value = 118.6364999999996

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 404251.78
$ws.Range("I32").Value = 4019.8645
$ws.Range("K32").Value = 4019.8645
$ws.Range("M32").Value = -3732.8645
$ws.Range("H54").Value = 10000
$ws.Range("J54").Value = 10000
$ws.Range("L54").Value = 10000
$ws.Range("N54").Value = -11538
$ws.Range("H61").Value = 1229.8125
$ws.Range("I61").Value = 1212.3
$ws.Range("J61").Value = 1259
$ws.Range("K61").Value = 1212.3
$ws.Range("L61").Value = 1259
$ws.Range("M61").Value = -1000.3
$ws.Range("N61").Value = -1683
$ws.Range("H74").Value = 1135.027
$ws.Range("I74").Value = 1060.1786
$ws.Range("J74").Value = 1367.8889
$ws.Range("K74").Value = 1060.1786
$ws.Range("L74").Value = 1367.8889
$ws.Range("M74").Value = -186.1786
$ws.Range("N74").Value = -3115.8889
$ws.Range("H77").Value = 1135.027
$ws.Range("I77").Value = 1060.1786
$ws.Range("J77").Value = 1367.8889
$ws.Range("K77").Value = 5300.893
$ws.Range("L77").Value = 6839.4445
$ws.Range("M77").Value = -932.893
$ws.Range("N77").Value = -15575.4445
$ws.Range("H132").Value = 16967606
$ws.Range("I132").Value = 20409388
$ws.Range("K132").Value = 61228164
$ws.Range("M132").Value = -61225634
$ws.Range("H136").Value = 1229.8125
$ws.Range("I136").Value = 1212.3
$ws.Range("J136").Value = 1259
$ws.Range("K136").Value = 3636.9
$ws.Range("L136").Value = 3777
$ws.Range("M136").Value = -1086.9
$ws.Range("N136").Value = -8877

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H32").Value = 4000
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").ClearContents()
$ws.Range("H107").Value = 1497.2
$ws.Range("I107").Value = 968.6667
$ws.Range("J107").Value = 2290
$ws.Range("K107").Value = 968.6667
$ws.Range("L107").Value = 2290
$ws.Range("M107").Value = 951.3333
$ws.Range("N107").Value = -6130
$ws.Range("H134").Value = 5995.75
$ws.Range("I134").Value = 2606.5454
$ws.Range("J134").Value = 11321.643
$ws.Range("K134").Value = 7819.6362
$ws.Range("L134").Value = 33964.929
$ws.Range("M134").Value = -5284.6362
$ws.Range("N134").Value = -39034.929

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 20842650
$ws.Range("I31").Value = 27781066
$ws.Range("J31").Value = 27400
$ws.Range("K31").Value = 27781066
$ws.Range("L31").Value = 27400
$ws.Range("M31").Value = -27780771
$ws.Range("N31").Value = -27990
$ws.Range("H34").Value = 20842650
$ws.Range("I34").Value = 27781066
$ws.Range("J34").Value = 27400
$ws.Range("K34").Value = 27781066
$ws.Range("L34").Value = 27400
$ws.Range("M34").Value = -27780864
$ws.Range("N34").Value = -27804
$ws.Range("H58").Value = 1102.6
$ws.Range("I58").Value = 1314
$ws.Range("J58").Value = 961.6667
$ws.Range("K58").Value = 1314
$ws.Range("L58").Value = 961.6667
$ws.Range("M58").Value = -1111
$ws.Range("N58").Value = -1367.6667
$ws.Range("H86").Value = 41633.5
$ws.Range("I86").Value = 1815.5454
$ws.Range("J86").Value = 90299.89
$ws.Range("K86").Value = 1815.5454
$ws.Range("L86").Value = 90299.89
$ws.Range("M86").Value = -692.5454
$ws.Range("N86").Value = -92545.89
$ws.Range("H89").Value = 41633.5
$ws.Range("I89").Value = 1815.5454
$ws.Range("J89").Value = 90299.89
$ws.Range("K89").Value = 9077.726999999999
$ws.Range("L89").Value = 451499.45
$ws.Range("M89").Value = -3461.726999999999
$ws.Range("N89").Value = -462731.45
$ws.Range("H132").Value = 44320.543
$ws.Range("I132").Value = 1760.8334
$ws.Range("K132").Value = 5282.5002
$ws.Range("M132").Value = -2752.5002
$ws.Range("H134").Value = 1877.8387
$ws.Range("I134").Value = 1434.0454
$ws.Range("J134").Value = 2962.6667
$ws.Range("K134").Value = 4302.1362
$ws.Range("L134").Value = 8888.000100000001
$ws.Range("M134").Value = -1767.1362
$ws.Range("N134").Value = -13958.0001
$ws.Range("H136").Value = 1102.6
$ws.Range("I136").Value = 1314
$ws.Range("J136").Value = 961.6667
$ws.Range("K136").Value = 3942
$ws.Range("L136").Value = 2885.0001
$ws.Range("M136").Value = -1392
$ws.Range("N136").Value = -7985.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1022.3333
$ws.Range("J34").Value = 1498.1818
$ws.Range("L34").Value = 4494.5454
$ws.Range("N34").Value = -4662.5454
$ws.Range("H113").Value = 980.3293
$ws.Range("I113").Value = 719.6667
$ws.Range("J113").Value = 1000.9079
$ws.Range("K113").Value = 2159.0001
$ws.Range("L113").Value = 3002.7237
$ws.Range("M113").Value = 10.9998999999998
$ws.Range("N113").Value = -7342.7237

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 424829.38
$ws.Range("I132").Value = 56662.445
$ws.Range("J132").Value = 1253205
$ws.Range("K132").Value = 169987.335
$ws.Range("L132").Value = 3759615
$ws.Range("M132").Value = -167457.335
$ws.Range("N132").Value = -3764675

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 39197.89
$ws.Range("I132").Value = 60908.47
$ws.Range("J132").Value = 2289.9
$ws.Range("K132").Value = 182725.41
$ws.Range("L132").Value = 6869.700000000001
$ws.Range("M132").Value = -180195.41
$ws.Range("N132").Value = -11929.7
$ws.Range("H136").Value = 12777
$ws.Range("I136").Value = 12146.167
$ws.Range("J136").Value = 15300.333
$ws.Range("K136").Value = 36438.501
$ws.Range("L136").Value = 45900.999
$ws.Range("M136").Value = -33888.501
$ws.Range("N136").Value = -51000.999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 3916.5
$ws.Range("I58").Value = 2874.75
$ws.Range("K58").Value = 2874.75
$ws.Range("M58").Value = -2566.75
$ws.Range("H132").Value = 63678350
$ws.Range("I132").Value = 102728890
$ws.Range("J132").Value = 2313218.8
$ws.Range("K132").Value = 308186670
$ws.Range("L132").Value = 6939656.399999999
$ws.Range("M132").Value = -308184140
$ws.Range("N132").Value = -6944716.399999999
